$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ B="Mapfre registra mais de 5.000 acionamentos de assistência no carnaval"; C="https://www.revistaapolice.com.br/2023/03/mapfre-registra-mais-de-5-000-acionamentos-de-assistencia-no-carnaval/"; D="Companhia reforçou atendimento na época para atender alta demanda. Assistências no Litoral Norte representaram 27% do total" },
    @{ B="Bradesco Saúde aprimora experiência do corretor com Inteligência Artificial"; C="https://www.revistaapolice.com.br/2023/03/bradesco-saude-aprimora-experiencia-do-corretor-com-inteligencia-artificial/"; D="Disponíveis no Portal de Negócios, novidades agilizam atendimento ao corretor no pós-venda dos produtos SPG" },
    @{ B="Allianz fecha parceria com o Rappi para distribuição de seguro via corretor"; C="https://www.revistaapolice.com.br/2023/03/allianz-fecha-parceria-com-o-rappi-para-distribuicao-de-seguro-via-corretor/"; D="Ação é voltada é exclusivamente aos clientes Rappi Prime e abrange todas as ofertas e pacotes do seguro Residencial da companhia" },
    @{ B="Mapfre anuncia Afonso Arinos como head Comercial de Agronegócio"; C="https://www.revistaapolice.com.br/2023/03/mapfre-apresenta-afonso-arinos-como-head-comercial-de-agronegocio/"; D="Com 20 anos de experiência em posições de liderança comercial, novo executivo irá liderar estratégias comerciais para a companhia" },
    @{ B="Akad disponibiliza seguro por assinatura para pequenas e médias empresas"; C="https://www.revistaapolice.com.br/2023/03/akad-disponibiliza-seguro-por-assinatura-para-pequenas-e-medias-empresas/"; D="Parceria com a MarketUP e a PDV Box é aposta da seguradora para chegar a mais de 150 mil usuários ativos" },
    @{ B="Tokio Marine registra crescimento de 40,6% na produção em 2022"; C="https://www.revistaapolice.com.br/2023/03/tokio-marine-registra-crescimento-de-406-na-producao-em-2022/"; D="A produção da companhia atingiu a marca de R$ 10,6 bilhões em faturamento, meta inicialmente prevista para ser alcançada apenas em 2024" },
    @{ B="PDMS quer levar seguro para fatia maior da população"; C="https://www.revistaapolice.com.br/2023/03/pdms-quer-levar-seguro-para-fatia-maior-da-populacao/"; D="Plano de Desenvolvimento do Mercado de Seguros possui 65 diretrizes em 4 áreas para aumentar o nível de proteção do País" },
    @{ B="Mapfre patrocina musical sobre a história do grupo Los Hermanos"; C="https://www.revistaapolice.com.br/2023/03/mapfre-patrocina-musical-sobre-a-historia-do-grupo-los-hermanos/"; D="Espetáculo, que tem a direção de Michel Melamed, estreia nesta sexta-feira, dia 17 de março, no Teatro Casa Grande, no Rio de Janeiro" },
    @{ B="SulAmérica lança linha de produtos Odonto Individual"; C="https://www.revistaapolice.com.br/2023/03/sulamerica-lanca-linha-de-produtos-odonto-individual/"; D="Desconto progressivo para dependentes e elegibilidade para toda a família estão entre os novos benefícios nos planos odontológicos da companhia" },
    @{ B="Generali apresenta resultado operacional de 6,5 bi de euros em 2022"; C="https://www.revistaapolice.com.br/2023/03/generali-apresenta-resultado-operacional-de-65-bi-de-euros-em-2022/"; D="De acordo com a companhia, o resultado operacional cresceu devido à expansão do segmento Vida, juntamente com Property and Casualty" }
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i].B
    $ws.Cells.Item($row, 3).Value = $data[$i].C
    $ws.Cells.Item($row, 4).Value = $data[$i].D
}
